$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Vcan", "Tlr2", "ECs", 2, 0.6666666666666666, 1.529781, 4.589343, 0.01315047351877542, 0.01315047351877542, 3, 1, 51.23401333333334, 153.70204, 0.9688226788583661, 0.9688226788583661, 78.37682015108, 705.39138135972, 0.01274047698271601, 0.01274047698271601),
    @("ECs", "Vcan", "Tlr2", "FAPs", 2, 0.6666666666666666, 1.529781, 4.589343, 0.01315047351877542, 0.01315047351877542, 2, 0.6666666666666666, 1.648742666666666, 4.946228, 0.03117732114163389, 0.03117732114163389, 2.522215205355999, 22.699936848204, 0.0004099965360594135, 0.0004099965360594136),
    @("FAPs", "Vcan", "Tlr2", "ECs", 3, 1, 103.676216, 311.028648, 0.8912330150752564, 0.8912330150752563, 3, 1, 51.23401333333334, 153.70204, 0.9688226788583661, 0.9688226788583661, 5311.748632893547, 47805.73769604192, 0.8634467571522285, 0.8634467571522284),
    @("FAPs", "Vcan", "Tlr2", "FAPs", 3, 1, 103.676216, 311.028648, 0.8912330150752564, 0.8912330150752563, 2, 0.6666666666666666, 1.648742666666666, 4.946228, 0.03117732114163389, 0.03117732114163389, 170.9354008377493, 1538.418607539744, 0.02778625792302791, 0.02778625792302791),
    @("sCs", "Vcan", "Tlr2", "ECs", 3, 1, 11.12297, 33.36891, 0.09561651140596822, 0.09561651140596822, 3, 1, 51.23401333333334, 153.70204, 0.9688226788583661, 0.9688226788583661, 569.8743932862667, 5128.869539576401, 0.09263544472342165, 0.09263544472342165),
    @("sCs", "Vcan", "Tlr2", "FAPs", 3, 1, 11.12297, 33.36891, 0.09561651140596822, 0.09561651140596822, 2, 0.6666666666666666, 1.648742666666666, 4.946228, 0.03117732114163389, 0.03117732114163389, 18.33891521905333, 165.05023697148, 0.002981066682546571, 0.002981066682546571)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $vals = $rows[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c+1).Value2 = $vals[$c]
    }
}

